$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.653.23'
$ws.Cells.Item(2, 5).Value = '  -0.21%  '

$ws.Cells.Item(3, 4).Value = '2.616.06'
$ws.Cells.Item(3, 5).Value = '  -0.43%  '

$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).Value = '594.09'
$ws.Cells.Item(5, 5).Value = '  -0.50%  '

$ws.Cells.Item(6, 4).Value = '150.68'
$ws.Cells.Item(6, 5).Value = '  +0.24%  '

$ws.Cells.Item(7, 5).Value = '  -0.03%  '

$ws.Cells.Item(8, 5).Value = '  -0.35%  '

$ws.Cells.Item(9, 4).Value = '0.113'
$ws.Cells.Item(9, 5).Value = '  +4.02%  '

$ws.Cells.Item(10, 2).Value = 'Toncoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(10, 4).Value = '5.80'
$ws.Cells.Item(10, 5).Value = '  +1.61%  '

$ws.Cells.Item(11, 2).Value = 'Cardano'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(11, 4).Value = '0.393'
$ws.Cells.Item(11, 5).Value = '  +2.99%  '

$ws.Cells.Item(12, 5).Value = '  +1.00%  '

$ws.Cells.Item(13, 4).Value = '27.81'
$ws.Cells.Item(13, 5).Value = '  +0.51%  '

$ws.Cells.Item(14, 4).Value = '3.084.98'
$ws.Cells.Item(14, 5).Value = '  -0.47%  '

$ws.Cells.Item(15, 4).Value = '63.567.53'
$ws.Cells.Item(15, 5).Value = '  -0.12%  '

$ws.Cells.Item(16, 4).Value = '0.0000168'
$ws.Cells.Item(16, 5).Value = '  +12.33%  '

$ws.Cells.Item(17, 4).Value = '2.590.02'
$ws.Cells.Item(17, 5).Value = '  -1.65%  '

$ws.Cells.Item(18, 4).Value = '12.14'
$ws.Cells.Item(18, 5).Value = '  -1.32%  '

$ws.Cells.Item(19, 4).Value = '4.76'
$ws.Cells.Item(19, 5).Value = '  +2.70%  '

$ws.Cells.Item(20, 4).Value = '346.80'
$ws.Cells.Item(20, 5).Value = '  -0.71%  '

$ws.Cells.Item(21, 4).Value = '6.97'
$ws.Cells.Item(21, 5).Value = '  +1.59%  '

$ws.Cells.Item(22, 4).Value = '1.00'
$ws.Cells.Item(22, 5).Value = '  +0.23%  '

$ws.Cells.Item(23, 5).Value = '  +1.47%  '

$ws.Cells.Item(24, 4).Value = '1.69'
$ws.Cells.Item(24, 5).Value = '  -2.57%  '

$ws.Cells.Item(25, 5).Value = '  +0.39%  '

$ws.Cells.Item(26, 4).Value = '9.13'
$ws.Cells.Item(26, 5).Value = '  -0.57%  '

$ws.Cells.Item(27, 4).Value = '8.29'
$ws.Cells.Item(27, 5).Value = '  +0.88%  '

$ws.Cells.Item(28, 4).Value = '546.44'
$ws.Cells.Item(28, 5).Value = '  -2.81%  '

$ws.Cells.Item(29, 4).Value = '0.161'
$ws.Cells.Item(29, 5).Value = '  -2.18%  '

$ws.Cells.Item(30, 4).Value = '1.00'
$ws.Cells.Item(30, 5).Value = '  -0.30%  '

$ws.Cells.Item(31, 4).Value = '0.0₃0900'
$ws.Cells.Item(31, 5).Value = '  +6.47%  '

$ws.Cells.Item(32, 5).Value = '  +0.39%  '

$ws.Cells.Item(33, 5).Value = '  +4.07%  '

$ws.Cells.Item(34, 5).Value = '  +2.38%  '

$ws.Cells.Item(35, 4).Value = '6.10'
$ws.Cells.Item(35, 5).Value = '  -0.14%  '

$ws.Cells.Item(36, 2).Value = 'Monero'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(36, 4).Value = '164.26'
$ws.Cells.Item(36, 5).Value = '  -2.84%  '

$ws.Cells.Item(37, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(37, 4).Value = '0.417'
$ws.Cells.Item(37, 5).Value = '  +1.96%  '

$ws.Cells.Item(38, 4).Value = '19.88'
$ws.Cells.Item(38, 5).Value = '  +2.70%  '

$ws.Cells.Item(39, 4).Value = '0.999'
$ws.Cells.Item(39, 5).Value = '  -0.08%  '

$ws.Cells.Item(40, 4).Value = '1.96'
$ws.Cells.Item(40, 5).Value = '  +0.88%  '

$ws.Cells.Item(41, 5).Value = '  -0.01%  '

$ws.Cells.Item(42, 4).Value = '167.15'
$ws.Cells.Item(42, 5).Value = '  -2.10%  '

$ws.Cells.Item(43, 4).Value = '4.08'
$ws.Cells.Item(43, 5).Value = '  +4.46%  '

$ws.Cells.Item(44, 4).Value = '23.16'
$ws.Cells.Item(44, 5).Value = '  +8.38%  '

$ws.Cells.Item(45, 4).Value = '0.0582'
$ws.Cells.Item(45, 5).Value = '  -2.54%  '

$ws.Cells.Item(46, 4).Value = '2.19'
$ws.Cells.Item(46, 5).Value = '  +9.94%  '

$ws.Cells.Item(47, 4).Value = '0.634'
$ws.Cells.Item(47, 5).Value = '  +0.54%  '

$ws.Cells.Item(48, 5).Value = '  +1.25%  '

$ws.Cells.Item(49, 5).Value = '  -0.21%  '

$ws.Cells.Item(50, 4).Value = '19.13'
$ws.Cells.Item(50, 5).Value = '  -0.44%  '

$ws.Cells.Item(51, 4).Value = '0.0₆0231'
$ws.Cells.Item(51, 5).Value = '  +18.04%  '
